$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: classical-best-embeddings vs. classical-best-tfidf ---
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.034
$ws.Range("E2").Value = 0.005
$ws.Range("F2").Value = 0.024
$ws.Range("G2").Value = 0.041
$ws.Range("H2").Value = 0.02
$ws.Range("I2").Value = 0.022
$ws.Range("J2").Value = 0.025

# --- Row 3: BERT-base vs. classical-best-tfidf ---
$ws.Range("C3").Value = 0.06
$ws.Range("D3").Value = 0.037
$ws.Range("E3").Value = 0.039
$ws.Range("F3").Value = 0.024
$ws.Range("G3").Value = 0.078
$ws.Range("H3").Value = 0.063
$ws.Range("I3").Value = 0.04
$ws.Range("J3").Value = 0.05

# --- Row 4: BERT-base vs. classical-best-embeddings ---
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = 0.026
$ws.Range("D4").Value = 0.014
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.037
$ws.Range("I4").Value = 0.019
$ws.Range("J4").Value = 0.026

# --- Row 5: BERT-base-nli vs. classical-best-tfidf ---
$ws.Range("B5").Value = 0.543
$ws.Range("C5").Value = 0.022
$ws.Range("D5").Value = 0.016
$ws.Range("E5").Value = 0.014
$ws.Range("F5").Value = -0.003
$ws.Range("G5").Value = 0.026
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.012
$ws.Range("J5").Value = 0.012

# --- Row 6: BERT-base-nli vs. classical-best-embeddings ---
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.543
$ws.Range("C6").Value = -0.012
$ws.Range("D6").Value = -0.007
$ws.Range("E6").Value = 0.008999999999999999
$ws.Range("F6").Value = -0.027
$ws.Range("G6").Value = -0.015
$ws.Range("H6").Value = -0.02
$ws.Range("I6").Value = -0.008999999999999999
$ws.Range("J6").Value = -0.012

# --- Row 7: BERT-base-nli vs. BERT-base ---
$ws.Range("B7").Value = 0.543
$ws.Range("C7").Value = -0.038
$ws.Range("D7").Value = -0.021
$ws.Range("E7").Value = -0.025
$ws.Range("G7").Value = -0.052
$ws.Range("H7").Value = -0.063
$ws.Range("I7").Value = -0.028
$ws.Range("J7").Value = -0.038
